$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update every cached "datetimeFigureOut" field (footer date placeholder)
#    from 07/09/2020 to 09/09/2020 across the slide master and every layout.
# ---------------------------------------------------------------------------
function Update-DateShape($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.HasTextFrame) {
            $tf = $sh.TextFrame
            if ($tf.HasText) {
                $txt = $tf.TextRange.Text
                if ($txt -eq "07/09/2020") {
                    $tf.TextRange.Text = "09/09/2020"
                }
            }
        }
    }
}

Update-DateShape $p.SlideMaster.Shapes
for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $lay = $p.SlideMaster.CustomLayouts.Item($i)
    Update-DateShape $lay.Shapes
}

# ---------------------------------------------------------------------------
# 2) On slide 7 (the Dashboard/flow diagram slide):
#    - move the "Notificações" rectangle to its new position
#    - duplicate it into a new "Slack" rectangle placed where the old
#      "Notificações" box used to sit
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)

# "Notificações" shape is named "Retângulo 40" in the deck.
$notifTargetName = "Ret" + [char]0x00E2 + "ngulo 40"
$notificacoes = $null
for ($i = 1; $i -le $s7.Shapes.Count; $i++) {
    $sh = $s7.Shapes.Item($i)
    if ($sh.Name -eq $notifTargetName) {
        $notificacoes = $sh
    }
}

# Duplicate the "Notificações" rectangle (copy/paste keeps the exact
# fill/line/style formatting) before moving the original, then rename,
# reposition and retext the copy into the new "Slack" rectangle.
$notificacoes.Copy()
$pasted = $s7.Shapes.Paste()
$slack = $pasted.Item(1)
$slack.Name = "Ret" + [char]0x00E2 + "ngulo 33"
$slack.Left = 381.36236220472443
$slack.Top = 450.31669691338584
$slack.TextFrame.TextRange.Text = "Slack"

# Now move the original "Notificações" rectangle to its new position.
$notificacoes.Left = 500.09260642519683
$notificacoes.Top = 407.19826771653544
